# "Organization of folders and cleaning of code"
#
# Observed changes (from the OOXML diff):
#   Sheet "Sheet1" (1st sheet):
#     - column A                width 14.1734693877551 -> 13.9030612244898
#     - columns B:C              width 33.2091836734694 -> 32.6683673469388
#     - columns D:AMK (rest)     width 14.1734693877551 -> 13.9030612244898
#   Sheet "Folha2" (2nd sheet):
#     - the view scrolled down / the active selection moved from B11 to A38
#       (topLeftCell A1 -> A46, selection B11 -> A38)
#     - column A                width 17.3622448979592 -> 17.1428571428571
#     - column B                width 31.5357142857143 -> 31.1836734693878
#     - columns C:AMK (rest)    width 11.5204081632653 -> 11.3418367346939
#
# The ColumnWidth COM property is expressed in "characters"; the engine
# converts it to the stored OOXML width with +5/6 on a 1/6-character grid,
# so the inputs below are the closest values on that grid to the widths
# seen in the target file.

$wb = $excel.ActiveWorkbook

# ----- Sheet1 : shrink the columns slightly -----
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Columns.Item(1).ColumnWidth = 13
$ws1.Range("B:C").ColumnWidth = 31.8333333333333
$ws1.Range("D:AMK").ColumnWidth = 13

# ----- Folha2 : scroll the window down and move the selection, and shrink columns -----
$ws2 = $wb.Worksheets.Item("Folha2")

$ws2.Columns.Item(1).ColumnWidth = 16.3333333333333
$ws2.Columns.Item(2).ColumnWidth = 30.3333333333333
$ws2.Range("C:AMK").ColumnWidth = 10.5

$win = $excel.ActiveWindow
$win.ScrollRow = 46
$win.ScrollColumn = 1
$ws2.Range("A38").Select()
